# Apply the change described in the diff:
#  - Rotate the I5:J18 data block down by one row (row 18's old values
#    wrap around to row 5), effectively reverting a previous shift.
#  - Update the active selection to I5:J5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Hat And No Cattle")

$firstRow = 5
$lastRow = 18

# Capture the current (pre-edit) values for I5:J18 so we can rotate them.
$oldValues = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $oldValues[$r] = @($ws.Cells.Item($r, 9).Value(), $ws.Cells.Item($r, 10).Value())
}

# Write rotated values: row r gets what used to be in row r-1,
# and the first row gets what used to be in the last row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    if ($r -eq $firstRow) {
        $src = $oldValues[$lastRow]
    } else {
        $src = $oldValues[$r - 1]
    }
    $ws.Cells.Item($r, 9).Value = $src[0]
    $ws.Cells.Item($r, 10).Value = $src[1]
}

# Update the selection shown in the saved sheet view.
$ws.Activate()
$ws.Range("I5:J5").Select()
